{"js": "// The empty \"ListParagraph\" paragraph (right after the \"odustao od praga\n// stjecanja, check box\" bullet) used to carry a stray `_GoBack` bookmark\n// left over from the last cursor position in the previous save. Remove it\n// from there ...\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ... and append a brand-new bulleted paragraph at the very end of the\n// body (same numbered list as the preceding items, numId 15) containing\n// the new sentence about searching/viewing/updating taxpayer data.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// insertParagraph clones the paragraph formatting (style, numbering,\n// justification, run formatting) of the paragraph it is inserted after,\n// matching how Word itself behaves when you press Enter at the end of a\n// list item.\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Unjeti podaci o poreznom obvezniku mogu se pretra\u017eivati, pregledavati i a\u017eurirati.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// The `_GoBack` bookmark now lands here instead, in the middle of the new\n// sentence (right after \"Unjeti poda\"), marking where the author's cursor\n// was when the document was saved.\nconst searchResults = newParagraph.search(\"Unjeti poda\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst splitPoint = searchResults.items[0].getRange(Word.RangeLocation.after);\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The stray \"_GoBack\" bookmark (left by the last saved cursor position)\n# used to sit in the empty ListParagraph right after the\n# \"odustao od praga stjecanja, check box\" bullet. Remove it from there.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n# Append a new bulleted paragraph at the very end of the document, in the\n# same numbered list (numId 15) as the paragraph before it.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($count)\n$fullText = \"Unjeti podaci o poreznom obvezniku mogu se pretra\u017eivati, pregledavati i a\u017eurirati.\"\n$newPara.Range.Text = $fullText\n\n# Re-insert the \"_GoBack\" bookmark in the middle of the new sentence,\n# right after \"Unjeti poda\", matching where the author's cursor was\n# when the document was last saved.\n$splitPos = $newPara.Range.Start + \"Unjeti poda\".Length\n$d.Bookmarks.Add(\"_GoBack\", $d.Range($splitPos, $splitPos))\n"}
